$wb = $excel.ActiveWorkbook

# Rename the second worksheet from "drug_categories" to "drug_shapes"
$ws2 = $wb.Worksheets.Item(2)
$ws2.Name = "drug_shapes"

# Make it the active sheet, with cell A12 selected (moves tabSelected
# from sheet1 "companies" to sheet2 "drug_shapes", and updates the
# workbook's activeTab / the sheet's selection accordingly)
$ws2.Activate()
$ws2.Range("A12").Select()
